# Adjusting model name capitalizations:
#   instructor-xl -> Instructor-XL
#   titan         -> Titan
# Only within the two "Discussion" paragraphs touched by the commit; the
# document has other, already-correctly-capitalized ("Titan", "Instructor-XL")
# and other-context lowercase ("instructor-xl" inside unrelated bullet/figure
# text) occurrences that must stay untouched, so every replacement below is
# anchored to a long, unique substring rather than a bare "titan"/"instructor-xl".

$d = $word.ActiveDocument

function Replace-Unique($searchText, $replaceText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($searchText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find.Execute did not find: $searchText"
    }
}

# Paragraph: "...combining the basic HyDE generator with either instructor-xl or titan as the embedding model..."
Replace-Unique "generator with either instructor-xl or titan as the embedding model" "generator with either Instructor-XL or Titan as the embedding model"

# Paragraph (figure caption, 10pt): "While both instructor-xl and titan embedding models showed initial promise..."
Replace-Unique "While both instructor-xl and titan embedding models showed initial promise" "While both Instructor-XL and Titan embedding models showed initial promise"

# "...The titan-generated embeddings achieved accuracies of 45% and 59%..."
Replace-Unique "The titan-generated embeddings achieved accuracies" "The Titan-generated embeddings achieved accuracies"

# "...In comparison, the instructor-generated embeddings demonstrated superior performance..."
Replace-Unique "In comparison, the instructor-generated embeddings demonstrated superior performance" "In comparison, the Instructor-generated embeddings demonstrated superior performance"

# "...we discontinued use of the titan embeddings and conducted all subsequent..."
Replace-Unique "we discontinued use of the titan embeddings" "we discontinued use of the Titan embeddings"

# "...using the instructor-xl model as the embedder."
Replace-Unique "evaluations using the instructor-xl model as the embedder" "evaluations using the Instructor-XL model as the embedder"

Write-Output "done"
